$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the sheet: "Заявки АХО" -> "Оперативная обстановка"
$ws.Name = "Оперативная обстановка"

# 2. Update the "Заявитель" header (column C) to "Инициатор / Заявитель"
$ws.Range("C1").Value = "Инициатор / Заявитель"

# 3. Add a new "Телефон" header in column I, matching the header row's
#    existing look & feel (copy format from the last header cell H1).
$ws.Range("I1").Value = "Телефон"
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)  # xlPasteFormats

# New column I should be the same width as the other header columns (15
# characters). ColumnWidth round-trips with a constant +5/6 offset in this
# engine, so compensate for it to land on exactly 15 in the saved file.
$ws.Columns.Item(9).ColumnWidth = 15 - 5/6
